# Update cryptos list data (prices and 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.024.62"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "2.298.07"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.26"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("E7").Value = "  -1.11%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.22"
$ws.Range("E10").Value = "  +7.48%  "

$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.19"
$ws.Range("E13").Value = "  +7.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.94"
$ws.Range("E14").Value = "  +1.88%  "

$ws.Range("D15").Value = "2.652.90"
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").Value = "2.311.77"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").Value = "42.924.54"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("E19").Value = "  +8.63%  "

$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.53"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("E24").Value = "  +9.59%  "

$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.95"
$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("E28").Value = "  +5.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.59"
$ws.Range("E29").Value = "  +1.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.18"
$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.02"
$ws.Range("E33").Value = "  +1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.67"
$ws.Range("E34").Value = "  +3.76%  "

$ws.Range("E35").Value = "  -1.34%  "

$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0690"
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("E38").Value = "  +1.64%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  -0.73%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("E41").Value = "  -0.42%  "

$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.30"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0291"
$ws.Range("E43").Value = "  +3.33%  "

$ws.Range("D44").Value = "1.975.00"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.19"
$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.55"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.57"
$ws.Range("E48").Value = "  +4.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.55"
$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("D50").Value = "2.520.41"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.80"
$ws.Range("E51").Value = "  +0.76%  "
